$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3351.1765
$ws.Range("I29").Value = 2800
$ws.Range("J29").Value = 3520.7693
$ws.Range("K29").Value = 8400
$ws.Range("L29").Value = 10562.3079
$ws.Range("M29").Value = -8119
$ws.Range("N29").Value = -11124.3079

$ws.Range("H33").Value = 1053.0526
$ws.Range("I33").Value = 175.45454
$ws.Range("K33").Value = 175.45454
$ws.Range("M33").Value = 53.54545999999999

$ws.Range("H75").Value = 31657
$ws.Range("J75").Value = 31657
$ws.Range("L75").Value = 31657
$ws.Range("N75").Value = -33529

$ws.Range("H78").Value = 31657
$ws.Range("J78").Value = 31657
$ws.Range("L78").Value = 94971
$ws.Range("N78").Value = -104331

$ws.Range("H93").Value = 28000
$ws.Range("J93").Value = 28000
$ws.Range("L93").Value = 28000
$ws.Range("N93").Value = -32992

$ws.Range("H132").Value = 590828.3
$ws.Range("I132").Value = 2754.6875
$ws.Range("J132").Value = 10000006
$ws.Range("K132").Value = 8264.0625
$ws.Range("L132").Value = 30000018
$ws.Range("M132").Value = -5734.0625
$ws.Range("N132").Value = -30005078

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H98").Value = 20327.5
$ws.Range("J98").Value = 20327.5
$ws.Range("L98").Value = 20327.5
$ws.Range("N98").Value = -26317.5

$ws.Range("H132").Value = 3618272.2
$ws.Range("I132").Value = 4251810
$ws.Range("K132").Value = 12755430
$ws.Range("M132").Value = -12752900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1229.4375
$ws.Range("I20").Value = 1266.2
$ws.Range("J20").Value = 1168.1666
$ws.Range("K20").Value = 1266.2
$ws.Range("L20").Value = 1168.1666
$ws.Range("M20").Value = -1019.2
$ws.Range("N20").Value = -1662.1666

$ws.Range("H94").Value = 1673.84
$ws.Range("I94").Value = 633.35297
$ws.Range("J94").Value = 3884.875
$ws.Range("K94").Value = 633.35297
$ws.Range("L94").Value = 3884.875
$ws.Range("M94").Value = -182.35297
$ws.Range("N94").Value = -4786.875

$ws.Range("H95").Value = 26250
$ws.Range("J95").Value = 26250
$ws.Range("L95").Value = 26250
$ws.Range("N95").Value = -31742

$ws.Range("H134").Value = 40891.73
$ws.Range("I134").Value = 1529.174
$ws.Range("J134").Value = 342671.34
$ws.Range("K134").Value = 4587.522
$ws.Range("L134").Value = 1028014.02
$ws.Range("M134").Value = -2052.522
$ws.Range("N134").Value = -1033084.02

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 24424.273
$ws.Range("I25").Value = 2270.3333
$ws.Range("J25").Value = 32732
$ws.Range("K25").Value = 2270.3333
$ws.Range("L25").Value = 32732
$ws.Range("M25").Value = -2096.3333
$ws.Range("N25").Value = -33080

$ws.Range("H96").Value = 20000
$ws.Range("J96").Value = 20000
$ws.Range("L96").Value = 20000
$ws.Range("N96").Value = -25492

$ws.Range("H122").Value = 1015.1111
$ws.Range("I122").Value = 856
$ws.Range("J122").Value = 1333.3334
$ws.Range("K122").Value = 2568
$ws.Range("L122").Value = 4000.0002
$ws.Range("M122").Value = -118
$ws.Range("N122").Value = -8900.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6251.6113
$ws.Range("I5").Value = 627.7143
$ws.Range("J5").Value = 9830.454
$ws.Range("K5").Value = 1883.1429
$ws.Range("L5").Value = 29491.362
$ws.Range("M5").Value = -1771.1429
$ws.Range("N5").Value = -29715.362

$ws.Range("H125").Value = 9147.666999999999
$ws.Range("J125").Value = 9147.666999999999
$ws.Range("L125").Value = 27443.001
$ws.Range("N125").Value = -37283.001

$ws.Range("H135").Value = 6251.6113
$ws.Range("I135").Value = 627.7143
$ws.Range("J135").Value = 9830.454
$ws.Range("K135").Value = 5649.428699999999
$ws.Range("L135").Value = 88474.086
$ws.Range("M135").Value = -3114.428699999999
$ws.Range("N135").Value = -93544.086

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0

$ws.Range("H63").Value = 29900
$ws.Range("J63").Value = 29900
$ws.Range("L63").Value = 29900
$ws.Range("N63").Value = -31272

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0

$ws.Range("H66").Value = 29900
$ws.Range("J66").Value = 29900
$ws.Range("L66").Value = 89700
$ws.Range("N66").Value = -96564

$ws.Range("H69").Value = 20000
$ws.Range("J69").Value = 20000
$ws.Range("L69").Value = 20000
$ws.Range("N69").Value = -21498

$ws.Range("H72").Value = 20000
$ws.Range("J72").Value = 20000
$ws.Range("L72").Value = 60000
$ws.Range("N72").Value = -67488

$ws.Range("H75").Value = 15000
$ws.Range("J75").Value = 15000
$ws.Range("L75").Value = 15000
$ws.Range("N75").Value = -16748

$ws.Range("H78").Value = 15000
$ws.Range("J78").Value = 15000
$ws.Range("L78").Value = 45000
$ws.Range("N78").Value = -53736

$ws.Range("H88").Value = 60000
$ws.Range("J88").Value = 60000
$ws.Range("L88").Value = 60000
$ws.Range("N88").Value = -60902

$ws.Range("H91").Value = 60000
$ws.Range("J91").Value = 60000
$ws.Range("L91").Value = 60000
$ws.Range("N91").Value = -63120

$ws.Range("H110").Value = 38000
$ws.Range("J110").Value = 38000
$ws.Range("L110").Value = 38000
$ws.Range("N110").Value = -46180

$ws.Range("H141").Value = 39379.6
$ws.Range("J141").Value = 39379.6
$ws.Range("L141").Value = 39379.6
$ws.Range("N141").Value = -49739.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 42178.5
$ws.Range("I13").Value = 1700
$ws.Range("J13").Value = 55671.332
$ws.Range("K13").Value = 1700
$ws.Range("L13").Value = 55671.332
$ws.Range("M13").Value = -1560
$ws.Range("N13").Value = -55951.332

$ws.Range("H40").Value = 48418.863
$ws.Range("I40").Value = 1954.5454
$ws.Range("K40").Value = 1954.5454
$ws.Range("M40").Value = -1818.5454

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = 0

$ws.Range("H95").Value = 10172
$ws.Range("J95").Value = 10172
$ws.Range("L95").Value = 10172
$ws.Range("N95").Value = -15664

$ws.Range("H97").Value = 27085.75
$ws.Range("J97").Value = 27085.75
$ws.Range("L97").Value = 27085.75
$ws.Range("N97").Value = -29067.75

$ws.Range("H132").Value = 753306.6
$ws.Range("I132").Value = 185018.55
$ws.Range("K132").Value = 555055.6499999999
$ws.Range("M132").Value = -552525.6499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1420.0358
$ws.Range("I81").Value = 850.7
$ws.Range("J81").Value = 1736.3334
$ws.Range("K81").Value = 1701.4
$ws.Range("L81").Value = 3472.6668
$ws.Range("M81").Value = -640.4000000000001
$ws.Range("N81").Value = -5594.6668

$ws.Range("H84").Value = 1420.0358
$ws.Range("I84").Value = 850.7
$ws.Range("J84").Value = 1736.3334
$ws.Range("K84").Value = 8507
$ws.Range("L84").Value = 17363.334
$ws.Range("M84").Value = -3203
$ws.Range("N84").Value = -27971.334

$ws.Range("H94").Value = 19714.285
$ws.Range("J94").Value = 18833.334
$ws.Range("L94").Value = 18833.334
$ws.Range("N94").Value = -20635.334

$ws.Range("H140").Value = 23995
$ws.Range("J140").Value = 23995
$ws.Range("L140").Value = 23995
$ws.Range("N140").Value = -34355

$ws.Range("H141").Value = 53232.5
$ws.Range("J141").Value = 53232.5
$ws.Range("L141").Value = 53232.5
$ws.Range("N141").Value = -63592.5
